# Version 2-8-17 Finalizado filtro Clases de Material
# Insert a new "CODIGO" column before the existing RUBRO column on the
# "Gasto Funcionamiento" sheet (sheet1), shifting B:F -> C:G, and add the
# header / merged cell for the new column B, matching the style used by
# the other header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift existing columns B:F one place to the right, opening up a blank
# column B for the new "CODIGO" field.
$ws.Columns.Item(2).Insert() | Out-Null

# The new column B9:B10 should be merged exactly like the other header
# cells (B9:B10 was the old B9:B10 merge before the shift, now living at
# C9:C10) - merge first, then clone the formatting from the header cell
# immediately to its right so the same cellXf (bold font + full border,
# centered) is reused instead of a new style being synthesized.
$ws.Range("B9:B10").Merge() | Out-Null

$ws.Range("C9").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null

$ws.Range("C10").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null

# Clear the clipboard marquee / keep the workbook tidy.
$excel.CutCopyMode = 0

# New header text for the inserted column.
$ws.Range("B9").Value = "CODIGO"

# Match the selection left behind by the editor after finishing the work.
$ws.Range("B11").Select() | Out-Null
